$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Januari 2019")

# Rename the first tab so it no longer references a specific month
$ws.Name = "Konsulttidrapport"

# Tighten the "hours worked" number format in column B (rows 8-38)
# from one decimal place to two decimal places
$ws.Range("B8:B38").NumberFormat = "0.00"

# Move the active selection as recorded in the saved view
$ws.Range("D46").Select() | Out-Null
